# Updates the cryptocurrency Price (column D) and Volume(1h) (column E)
# values for rows 2-51 of the active worksheet, matching the latest
# scrape performed by the GitHub Actions workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, new Price (D) text, new Volume(1h) (E) text. A value of $null means
# that column is unchanged for that row.
$updates = @(
    @{ Row = 2; D = "24.746.56"; E = "  +0.80%  " },
    @{ Row = 3; D = "1.700.12"; E = "  +0.40%  " },
    @{ Row = 4; D = "1.003"; E = "  +0.20%  " },
    @{ Row = 5; D = "316.45"; E = "  +0.43%  " },
    @{ Row = 6; D = $null; E = "  +0.18%  " },
    @{ Row = 7; D = "0.3942"; E = "  +0.27%  " },
    @{ Row = 8; D = "0.4064"; E = "  +1.76%  " },
    @{ Row = 9; D = "1.512"; E = "  -0.63%  " },
    @{ Row = 10; D = "1.004"; E = "  +0.36%  " },
    @{ Row = 11; D = "53.20"; E = "  +2.03%  " },
    @{ Row = 12; D = "0.08897"; E = "  +1.93%  " },
    @{ Row = 13; D = "7.662"; E = "  +6.45%  " },
    @{ Row = 14; D = "23.75"; E = "  +3.02%  " },
    @{ Row = 15; D = "8.175"; E = "  +7.26%  " },
    @{ Row = 16; D = "0.00001327"; E = "  +1.22%  " },
    @{ Row = 17; D = "1.691.55"; E = "  +0.19%  " },
    @{ Row = 18; D = "99.61"; E = "  +0.17%  " },
    @{ Row = 19; D = "0.07113"; E = "  +0.78%  " },
    @{ Row = 20; D = "19.88"; E = "  +1.66%  " },
    @{ Row = 21; D = "7.141"; E = "  +3.88%  " },
    @{ Row = 22; D = "1.005"; E = "  +0.47%  " },
    @{ Row = 23; D = "14.73"; E = "  +5.02%  " },
    @{ Row = 24; D = "24.717.99"; E = "  +0.76%  " },
    @{ Row = 25; D = "3.160"; E = "  +3.30%  " },
    @{ Row = 26; D = "2.352"; E = "  +0.51%  " },
    @{ Row = 27; D = "23.07"; E = "  +3.72%  " },
    @{ Row = 28; D = "9.184"; E = "  +21.58%  " },
    @{ Row = 29; D = "164.44"; E = "  +1.74%  " },
    @{ Row = 30; D = "139.23"; E = "  +3.88%  " },
    @{ Row = 31; D = "5.144"; E = "  -1.40%  " },
    @{ Row = 32; D = $null; E = "  +8.79%  " },
    @{ Row = 33; D = "0.09031"; E = "  +5.91%  " },
    @{ Row = 34; D = "1.072"; E = "  -0.99%  " },
    @{ Row = 35; D = "0.03027"; E = "  +10.44%  " },
    @{ Row = 36; D = "0.2787"; E = "  +2.89%  " },
    @{ Row = 37; D = "11.13"; E = "  -1.04%  " },
    @{ Row = 38; D = "1.959"; E = "  +1.39%  " },
    @{ Row = 39; D = "14.44"; E = "  +0.25%  " },
    @{ Row = 40; D = "0.09271"; E = "  +2.64%  " },
    @{ Row = 41; D = $null; E = "  +1.61%  " },
    @{ Row = 42; D = "1.466"; E = "  -0.67%  " },
    @{ Row = 43; D = "16.18"; E = "  +4.93%  " },
    @{ Row = 44; D = "2.631"; E = "  +4.57%  " },
    @{ Row = 45; D = "0.7234"; E = "  +1.22%  " },
    @{ Row = 46; D = "4.225"; E = "  +0.48%  " },
    @{ Row = 47; D = "1.362"; E = "  +2.45%  " },
    @{ Row = 48; D = "1.002"; E = "  +0.17%  " },
    @{ Row = 49; D = "140.30"; E = "  -0.10%  " },
    @{ Row = 50; D = "0.07994"; E = "  +0.06%  " },
    @{ Row = 51; D = "89.94"; E = "  +2.36%  " }
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.D -ne $null) {
        $cell = $ws.Cells.Item($row, 4)
        # Many of the price strings (e.g. "1.003") look like plain numbers to
        # Excel and would otherwise be silently re-typed as numeric values,
        # losing the original text formatting (e.g. "53.20" -> 53.2). Force
        # the cell to Text before assigning, then restore the default
        # "Normal" style so no visible formatting change is introduced.
        if ($item.D -match '^-?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
            $cell.Value = $item.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $item.D
        }
    }

    if ($item.E -ne $null) {
        $ws.Cells.Item($row, 5).Value = $item.E
    }
}
